$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.980.74'
$ws.Range('E2').Value = '  -2.36%  '
$ws.Range('D3').Value = '1.963.78'
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.44'
$ws.Range('E5').Value = '  -1.30%  '
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4966'
$ws.Range('E7').Value = '  -0.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4199'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.85'
$ws.Range('E9').Value = '  -2.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09185'
$ws.Range('E10').Value = '  +3.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.096'
$ws.Range('E11').Value = '  -2.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.76'
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').Value = '1.963.83'
$ws.Range('E13').Value = '  -15.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.869'
$ws.Range('E14').Value = '  -3.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.452'
$ws.Range('E15').Value = '  -1.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.009'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.59'
$ws.Range('E17').Value = '  -5.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001099'
$ws.Range('E18').Value = '  -1.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06719'
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.25'
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.006'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.942'
$ws.Range('E22').Value = '  -1.55%  '
$ws.Range('D23').Value = '29.036.59'
$ws.Range('E23').Value = '  -2.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.99'
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.264'
$ws.Range('E25').Value = '  -1.83%  '
$ws.Range('D26').Value = '2.198.08'
$ws.Range('E26').Value = '  -10.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.61'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.58'
$ws.Range('E28').Value = '  -1.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.365'
$ws.Range('E29').Value = '  -2.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.249'
$ws.Range('E30').Value = '  -4.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '126.51'
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.043'
$ws.Range('E32').Value = '  -1.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09834'
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.515'
$ws.Range('E34').Value = '  -2.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.811'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.713'
$ws.Range('E36').Value = '  -3.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02430'
$ws.Range('E37').Value = '  -1.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.318'
$ws.Range('E38').Value = '  +1.86%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06366'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.001'
$ws.Range('E40').Value = '  -7.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6443'
$ws.Range('E41').Value = '  -1.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.40'
$ws.Range('E42').Value = '  -4.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1983'
$ws.Range('E43').Value = '  -4.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.006'
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6205'
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.338'
$ws.Range('E46').Value = '  +5.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.194'
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '13.27'
$ws.Range('E48').Value = '  -1.50%  '
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000324'
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06972'
$ws.Range('E51').Value = '  -0.96%  '
